# MyTestAccountPasswords.xlsx update
# - minimize the workbook window
# - add a new header row (id / username / pasword) above the existing data
# - add 3 new shared strings used by the header row
# - re-point the hyperlinks to follow the shifted rows (and re-shuffle a couple
#   of the hyperlinks lower on the sheet, matching the authoritative edit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. minimize the workbook window -----------------------------------
$win = $excel.ActiveWindow
$win.WindowState = -4140

# --- 2. push the first block of rows (1-9) down by one row -------------
# Copy bottom-up so we never clobber a row before it has been copied.
for ($r = 9; $r -ge 1; $r--) {
    $srcRow = $ws.Range("A" + $r + ":B" + $r)
    $dstCell = $ws.Range("A" + ($r + 1))
    $srcRow.Copy($dstCell)
}

# --- 3. write the new header row (row 1) --------------------------------
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "pasword"
$ws.Range("C1").Value = "id"

# --- 4. rebuild the hyperlinks collection -------------------------------
# Drop every existing hyperlink; they will be re-added below in the exact
# order needed so the relationship ids line up the same way the
# authoritative edit produced them.
$guard = 0
while ($ws.Hyperlinks.Count -gt 0 -and $guard -lt 50) {
    $current = @($ws.Hyperlinks)
    $current[0].Delete()
    $guard = $guard + 1
}

# rId1..rId9 for the shifted rows, rId10 for A20, rId11 for A18, rId12 for A17.
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:mfoibrh_seligsteinbergescuskysensteinsonmanwitz_1372209828@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:gyskysm_schrockescu_1372209832@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:ybsmoen_thurnberg_1372209830@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:gevhsig_schrockwitz_1372209825@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:etlcvqj_carrierostein_1372209822@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:tnhqrvn_greeneescu_1372209820@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:efkqolf_rosenthalescu_1372209818@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:yiasuby_bowersman_1372209815@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:wfvkakm_moidusky_1372209813@tfbnw.net") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A20"), "mailto:screensaver_catnvpn_screensaver@tfbnw.net") | Out-Null

# A18 keeps its legacy "display" text that differs from the cell's actual
# text (mailto:-prefixed). Setting TextToDisplay also rewrites the cell
# value, so set it explicitly and then restore the original cell text.
$h18 = $ws.Hyperlinks.Add($ws.Range("A18"), "mailto:savingsonscreen_xkbqzhl_savingsonscreen@tfbnw.net", [Type]::Missing, [Type]::Missing, "mailto:savingsonscreen_xkbqzhl_savingsonscreen@tfbnw.net")
$ws.Range("A18").Value = "savingsonscreen_xkbqzhl_savingsonscreen@tfbnw.net"

$ws.Hyperlinks.Add($ws.Range("A17"), "mailto:screensavingsapp@gmail.com") | Out-Null

# --- 5. fix up the sheet selection --------------------------------------
$ws.Range("B17").Select()
